$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds the last-changed date as an Excel serial
# number. This automatic update bumps every existing 45180 (2023-09-11)
# value to 45181 (2023-09-12) for rows 2 through 135.
$range = $ws.Range("C2:C135")
foreach ($cell in $range.Cells) {
    if ($cell.Value2 -eq 45180) {
        $cell.Value2 = 45181
    }
}
